# Added the stuff at the bottom of the sheets.
# - Practice rows (2-5) get a "pair_kind" value ("generic") in column J, and the
#   carrier word moves from column C into K (C/D above already held "practice").
# - Practice-stimulus detail rows (6-13) get their "kind" column switched to "generic".
# - A brand new block is appended starting at row 27: a "stim details" header,
#   a table header row (month/word_type/need_audio/need_image/word/count/find images),
#   and 8 data rows describing video/audio counts for months 6 and 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'p1'
$ws.Range("C2").Value = 'practice'
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 'A'
$ws.Range("I2").Value = 'banana_kitty'
$ws.Range("J2").Value = 'generic'
$ws.Range("K2").Value = 'can'

$ws.Range("A3").Value = 'p2'
$ws.Range("C3").Value = 'practice'
$ws.Range("H3").Value = 'B'
$ws.Range("I3").Value = 'bear_cracker'
$ws.Range("J3").Value = 'generic'
$ws.Range("K3").Value = 'do'

$ws.Range("A4").Value = 'p3'
$ws.Range("C4").Value = 'practice'
$ws.Range("H4").Value = 'C'
$ws.Range("I4").Value = 'hair_cup'
$ws.Range("J4").Value = 'generic'
$ws.Range("K4").Value = 'look'

$ws.Range("A5").Value = 'p4'
$ws.Range("C5").Value = 'practice'
$ws.Range("H5").Value = 'D'
$ws.Range("I5").Value = 'cheerios_water'
$ws.Range("J5").Value = 'generic'
$ws.Range("K5").Value = 'where'

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 'banana'
$ws.Range("C6").Value = 'generic'
$ws.Range("D6").Value = 'can'
$ws.Range("H6").Value = 'E'
$ws.Range("K6").Value = 'can'

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 'kitty'
$ws.Range("C7").Value = 'generic'
$ws.Range("D7").Value = 'can'
$ws.Range("H7").Value = 'F'
$ws.Range("K7").Value = 'do'

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 'bear'
$ws.Range("C8").Value = 'generic'
$ws.Range("D8").Value = 'do'
$ws.Range("H8").Value = 'G'
$ws.Range("K8").Value = 'look'

$ws.Range("A9").Value = 4
$ws.Range("B9").Value = 'cracker'
$ws.Range("C9").Value = 'generic'
$ws.Range("D9").Value = 'do'
$ws.Range("H9").Value = 'H'
$ws.Range("K9").Value = 'where'

$ws.Range("A10").Value = 5
$ws.Range("B10").Value = 'cup'
$ws.Range("C10").Value = 'generic'
$ws.Range("D10").Value = 'look'

$ws.Range("A11").Value = 6
$ws.Range("B11").Value = 'hair'
$ws.Range("C11").Value = 'generic'
$ws.Range("D11").Value = 'look'

$ws.Range("A12").Value = 7
$ws.Range("B12").Value = 'cheerios'
$ws.Range("C12").Value = 'generic'
$ws.Range("D12").Value = 'where'

$ws.Range("A13").Value = 8
$ws.Range("B13").Value = 'water'
$ws.Range("C13").Value = 'generic'
$ws.Range("D13").Value = 'where'

$ws.Range("A27").Value = 'stim details'

$ws.Range("A28").Value = 'month'
$ws.Range("B28").Value = 'word_type'
$ws.Range("C28").Value = 'need_audio'
$ws.Range("D28").Value = 'need_image'
$ws.Range("E28").Value = 'word'
$ws.Range("F28").Value = 'count'
$ws.Range("G28").Value = 'find images'

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = 'video'

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = 'video'

$ws.Range("A31").Value = 7
$ws.Range("B31").Value = 'video'

$ws.Range("A32").Value = 7
$ws.Range("B32").Value = 'video'

$ws.Range("A33").Value = 6
$ws.Range("B33").Value = 'audio'

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = 'audio'

$ws.Range("A35").Value = 7
$ws.Range("B35").Value = 'audio'

$ws.Range("A36").Value = 7
$ws.Range("B36").Value = 'audio'
